$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 into the two new header cells, then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I ("I0") and J ("IF"), rows 2-14
$values = @(
    @(8,8),
    @(8,9),
    @(8,8),
    @(5,5),
    @(4,4),
    @(8,8),
    @(9,9),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(5,5),
    @(7,7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
